# Update cryptocurrency price/volume data on cryptos.xlsx Sheet1
# Mirrors the GitHub Actions crypto-data refresh commit (Sun Jul 14 21:11:29 UTC 2024).
# All target cells (Price / Volume(1h) columns, plus the Aptos/Fetch.AI name+link
# swap in rows 33-34) are stored as plain text, not numbers. Each new literal is
# written with a leading apostrophe so Excel treats it as text instead of
# auto-converting it to a number/date; the Style is then reset to "Normal" so the
# apostrophe's "quote prefix" marker doesn't leave a stray style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.222.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +2.78%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.206.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +1.38%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.14%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''539.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.90%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''146.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +4.77%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.05%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -2.10%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  +0.92%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D11").Value = '''0.433'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -0.83%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''3.750.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +1.12%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -1.88%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''25.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.67%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  +1.16%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''60.179.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +2.60%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''3.188.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +1.51%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''6.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.61%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''13.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +2.46%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''8.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +1.66%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''371.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -1.19%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.19%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  -1.11%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''69.61'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -0.21%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +1.51%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''8.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +4.85%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.06%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''0.0₃0879'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +1.91%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''22.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.32%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  +0.81%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''6.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +1.66%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +2.72%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = '''Fetch.AI'
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = '''1.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +3.60%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = '''Aptos'
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = '''6.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +4.36%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''158.83'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +1.40%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''1.38'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +3.52%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''26.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +4.98%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.797.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +4.43%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.0315'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +9.28%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E41").Value = '''  +0.70%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''4.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -1.38%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''39.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +2.13%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  -0.17%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.106'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +1.63%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''3.241.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +1.10%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.986'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +0.83%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''6.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.92%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''20.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +3.54%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.798'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +6.74%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -0.13%  '
$ws.Range("E51").Style = "Normal"
